{"js": "// Append new \"Sonar\" bullet content (and sub-bullets) at the end of the\n// document's Notes list, reusing the existing trailing empty list\n// paragraph for the \"Sonar\" heading line, then adding the detail bullets\n// below it at the appropriate outline levels.\n\nconst body = context.document.body;\nconst lastParagraph = body.paragraphs.getLast();\nawait context.sync();\n\n// The trailing paragraph is currently empty (still has the ListParagraph /\n// numId=1 / ilvl=0 list formatting) -- give it the \"Sonar\" text.\nlastParagraph.insertText(\"Sonar\", \"End\");\n\n// Helper: insert a new list paragraph after `anchor`, set its outline\n// level, and return it so later calls can keep appending after it.\nfunction addBullet(anchor, text, level) {\n  const p = anchor.insertParagraph(text, \"After\");\n  p.listItem.level = level;\n  return p;\n}\n\nlet cursor = lastParagraph;\n\ncursor = addBullet(\n  cursor,\n  \"Technical debt = estimated time required to fix all the problems\",\n  1\n);\n\ncursor = addBullet(\n  cursor,\n  \"Code smell = characteristic in the cod that indicates a potential problem\",\n  1\n);\n\ncursor = addBullet(cursor, \"What \", 2);\ncursor.getRange(\"End\").insertText(\"constitutes a \", \"End\");\ncursor.getRange(\"End\").insertText(\"potential problem is subjective\", \"End\");\n\ncursor = addBullet(\n  cursor,\n  \"Usually indicates issues with long-term maintainability of the code \\u2013 it can be improved \",\n  2\n);\ncursor\n  .getRange(\"End\")\n  .insertText(\n    \"to ensure future developers have a minimal risk of creating an error DUE TO THAT CODE\",\n    \"End\"\n  );\n\ncursor = addBullet(\n  cursor,\n  \"Duplication = locations where identical code is \",\n  1\n);\ncursor.getRange(\"End\").insertText(\"located & ought to be their own methods\", \"End\");\n\ncursor = addBullet(cursor, \"Quality gates and Fixing the water leaks\", 1);\n\nawait context.sync();\n", "ps1": "# Append the new \"Sonar\" bullet (and its sub-bullets) to the end of the\n# Notes list. The document already ends with an empty ListParagraph /\n# numId=1 / level-1 bullet, so it becomes the \"Sonar\" line; the detail\n# bullets are then added below it at the correct outline levels.\n# Note: Word COM's ListFormat.ListLevelNumber is 1-based (1 = ilvl 0).\n\n$d = $word.ActiveDocument\n\n$last = $d.Paragraphs.Last\n$last.Range.InsertAfter(\"Sonar\")\n\n$last.Range.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.InsertAfter(\"Technical debt = estimated time required to fix all the problems\")\n$p.Range.ListFormat.ListLevelNumber = 2\n\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.InsertAfter(\"Code smell = characteristic in the cod that indicates a potential problem\")\n$p.Range.ListFormat.ListLevelNumber = 2\n\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.InsertAfter(\"What \")\n$p.Range.ListFormat.ListLevelNumber = 3\n$r = $p.Range\n$r.Collapse(0)\n$r.InsertAfter(\"constitutes a \")\n$r = $p.Range\n$r.Collapse(0)\n$r.InsertAfter(\"potential problem is subjective\")\n\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.InsertAfter(\"Usually indicates issues with long-term maintainability of the code \u2013 it can be improved \")\n$p.Range.ListFormat.ListLevelNumber = 3\n$r = $p.Range\n$r.Collapse(0)\n$r.InsertAfter(\"to ensure future developers have a minimal risk of creating an error DUE TO THAT CODE\")\n\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.InsertAfter(\"Duplication = locations where identical code is \")\n$p.Range.ListFormat.ListLevelNumber = 2\n$r = $p.Range\n$r.Collapse(0)\n$r.InsertAfter(\"located & ought to be their own methods\")\n\n$p.Range.InsertParagraphAfter()\n$p = $d.Paragraphs.Last\n$p.Range.InsertAfter(\"Quality gates and Fixing the water leaks\")\n$p.Range.ListFormat.ListLevelNumber = 2\n"}
